$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9613366723060608
$ws.Range("B1").Value = 1.585789322853088
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.619231700897217
$ws.Range("E1").Value = 1.354471564292908
